# Applies the "Deploying to gh-pages" metadata refresh:
#   1. Rename the "Include from SerumTumorMarker" sheet to "Include #0"
#   2. Insert a new "Jurisdiction" metadata row (with an empty value) right
#      after the existing "Contact" row on the Metadata sheet
#   3. Bump the "Date" metadata value to the new publish timestamp

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$include = $wb.Worksheets.Item(2)

# 1. Rename the second sheet.
$include.Name = "Include #0"

# 2. Make room for the new "Jurisdiction" row right below "Contact" (row 10),
#    pushing Description/Purpose/Copyright/Immutable down by one row.
$meta.Rows.Item(11).Insert()

# Carry over the bordered/wrapped "data row" look from the row above instead
# of leaving the blank style Excel assigns new rows by default.
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# 3. Refresh the publish Date value.
$meta.Range("B8").Value = "2024-09-17T19:55:11+00:00"
